$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new columns at D:E, shifting existing D:K data to F:M
$ws.Columns("D:E").Insert()

# Step 2: Copy formatting (number format/font/style) from column F (already correctly
# styled per-row) into the new D:E columns so the new cells match their row style.
$ws.Range("F5:G102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Step 3: Write final literal values for every data cell in columns D:M, row by row,
# reflecting the refreshed quarterly financial data (two new quarters added at D:E,
# and the scraped historical quarters re-stated per the source).
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("F7").Value = 43281
$ws.Range("G7").Value = 43190
$ws.Range("H7").Value = 43100
$ws.Range("I7").Value = 43008
$ws.Range("J7").Value = 42916
$ws.Range("K7").Value = 42825
$ws.Range("L7").Value = 42735
$ws.Range("M7").Value = 42643
$ws.Range("D8").Value = 941200
$ws.Range("E8").Value = 1104800
$ws.Range("F8").Value = 924500
$ws.Range("G8").Value = 858500
$ws.Range("H8").Value = 1002600
$ws.Range("I8").Value = 934200
$ws.Range("J8").Value = 813500
$ws.Range("K8").Value = 817900
$ws.Range("L8").Value = 817900
$ws.Range("M8").Value = 722900
$ws.Range("D9").Value = 733200
$ws.Range("E9").Value = 887000
$ws.Range("F9").Value = 697800
$ws.Range("G9").Value = 630700
$ws.Range("H9").Value = 802800
$ws.Range("I9").Value = 729000
$ws.Range("J9").Value = 627300
$ws.Range("K9").Value = 556700
$ws.Range("L9").Value = 776200
$ws.Range("M9").Value = 596400
$ws.Range("D10").Value = 208000
$ws.Range("E10").Value = 217800
$ws.Range("F10").Value = 226600
$ws.Range("G10").Value = 227800
$ws.Range("H10").Value = 199800
$ws.Range("I10").Value = 205100
$ws.Range("J10").Value = 186200
$ws.Range("K10").Value = 261200
$ws.Range("L10").Value = 41800
$ws.Range("M10").Value = 126500
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "NA"
$ws.Range("H12").Value = "NA"
$ws.Range("I12").Value = "NA"
$ws.Range("J12").Value = "NA"
$ws.Range("K12").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("D15").Value = 2700
$ws.Range("E15").Value = 2600
$ws.Range("F15").Value = 2500
$ws.Range("G15").Value = 2400
$ws.Range("H15").Value = 2500
$ws.Range("I15").Value = 2800
$ws.Range("J15").Value = 3400
$ws.Range("K15").Value = 3000
$ws.Range("L15").Value = 2800
$ws.Range("M15").Value = 2800
$ws.Range("D17").Value = 797100
$ws.Range("E17").Value = 951400
$ws.Range("F17").Value = 759500
$ws.Range("G17").Value = 707100
$ws.Range("H17").Value = 914100
$ws.Range("I17").Value = 817300
$ws.Range("J17").Value = 679400
$ws.Range("K17").Value = 613800
$ws.Range("L17").Value = 830600
$ws.Range("M17").Value = 659600
$ws.Range("D18").Value = 144100
$ws.Range("E18").Value = 153400
$ws.Range("F18").Value = 164900
$ws.Range("G18").Value = 151400
$ws.Range("H18").Value = 88500
$ws.Range("I18").Value = 116900
$ws.Range("J18").Value = 134100
$ws.Range("K18").Value = 204200
$ws.Range("L18").Value = -12700
$ws.Range("M18").Value = 63300
$ws.Range("D20").Value = 22200
$ws.Range("E20").Value = 21900
$ws.Range("F20").Value = 30200
$ws.Range("G20").Value = 36700
$ws.Range("H20").Value = -2200
$ws.Range("I20").Value = 35700
$ws.Range("J20").Value = -1900
$ws.Range("K20").Value = 28800
$ws.Range("L20").Value = -18800
$ws.Range("M20").Value = 17100
$ws.Range("D21").Value = 216200
$ws.Range("E21").Value = 223300
$ws.Range("F21").Value = 243900
$ws.Range("G21").Value = 233600
$ws.Range("H21").Value = 133100
$ws.Range("I21").Value = 199300
$ws.Range("J21").Value = 179400
$ws.Range("K21").Value = 278400
$ws.Range("L21").Value = 12300
$ws.Range("M21").Value = 124900
$ws.Range("D22").Value = 57500
$ws.Range("E22").Value = 54200
$ws.Range("F22").Value = 57100
$ws.Range("G22").Value = 54600
$ws.Range("H22").Value = 58400
$ws.Range("I22").Value = 61300
$ws.Range("J22").Value = 65400
$ws.Range("K22").Value = 67500
$ws.Range("L22").Value = 71300
$ws.Range("M22").Value = 71600
$ws.Range("D23").Value = 108800
$ws.Range("E23").Value = 121100
$ws.Range("F23").Value = 138000
$ws.Range("G23").Value = 133600
$ws.Range("H23").Value = 28000
$ws.Range("I23").Value = 91300
$ws.Range("J23").Value = 66900
$ws.Range("K23").Value = 165500
$ws.Range("L23").Value = -102700
$ws.Range("M23").Value = 8900
$ws.Range("D24").Value = 8600
$ws.Range("E24").Value = 28500
$ws.Range("F24").Value = 47600
$ws.Range("G24").Value = 46500
$ws.Range("H24").Value = -6200
$ws.Range("I24").Value = -15600
$ws.Range("J24").Value = 28200
$ws.Range("K24").Value = 61900
$ws.Range("L24").Value = -36300
$ws.Range("M24").Value = 27500
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("D26").Value = 100200
$ws.Range("E26").Value = 92600
$ws.Range("F26").Value = 90400
$ws.Range("G26").Value = 87100
$ws.Range("H26").Value = 34200
$ws.Range("I26").Value = 106800
$ws.Range("J26").Value = 38700
$ws.Range("K26").Value = 103500
$ws.Range("L26").Value = -66400
$ws.Range("M26").Value = -18600
$ws.Range("D27").Value = 95600
$ws.Range("E27").Value = 91000
$ws.Range("F27").Value = 89000
$ws.Range("G27").Value = 85200
$ws.Range("H27").Value = 19000
$ws.Range("I27").Value = 104900
$ws.Range("J27").Value = 35900
$ws.Range("K27").Value = 101800
$ws.Range("L27").Value = -52400
$ws.Range("M27").Value = -21600
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("D32").Value = -22200
$ws.Range("E32").Value = -21900
$ws.Range("F32").Value = -30200
$ws.Range("G32").Value = -36700
$ws.Range("H32").Value = 2200
$ws.Range("I32").Value = -35700
$ws.Range("J32").Value = 1900
$ws.Range("K32").Value = -28800
$ws.Range("L32").Value = 18800
$ws.Range("M32").Value = -17100
$ws.Range("D33").Value = 95600
$ws.Range("E33").Value = 91000
$ws.Range("F33").Value = 89000
$ws.Range("G33").Value = 85200
$ws.Range("H33").Value = 19000
$ws.Range("I33").Value = 104900
$ws.Range("J33").Value = 35900
$ws.Range("K33").Value = 101800
$ws.Range("L33").Value = -52400
$ws.Range("M33").Value = -21600
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("D35").Value = 95600
$ws.Range("E35").Value = 91000
$ws.Range("F35").Value = 89000
$ws.Range("G35").Value = 85200
$ws.Range("H35").Value = 19000
$ws.Range("I35").Value = 104900
$ws.Range("J35").Value = 35900
$ws.Range("K35").Value = 101800
$ws.Range("L35").Value = -52400
$ws.Range("M35").Value = -21600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("F38").Value = 43281
$ws.Range("G38").Value = 43190
$ws.Range("H38").Value = 43100
$ws.Range("I38").Value = 43008
$ws.Range("J38").Value = 42916
$ws.Range("K38").Value = 42825
$ws.Range("L38").Value = 42735
$ws.Range("M38").Value = 42643
$ws.Range("D41").Value = 499500
$ws.Range("E41").Value = 219900
$ws.Range("F41").Value = 436900
$ws.Range("G41").Value = 435200
$ws.Range("H41").Value = 266700
$ws.Range("I41").Value = 334700
$ws.Range("J41").Value = 365100
$ws.Range("K41").Value = 234600
$ws.Range("L41").Value = 243600
$ws.Range("M41").Value = 351700
$ws.Range("D42").Value = 140000
$ws.Range("E42").Value = 105000
$ws.Range("F42").Value = 85900
$ws.Range("G42").Value = 84200
$ws.Range("H42").Value = 44300
$ws.Range("I42").Value = 38300
$ws.Range("J42").Value = 14700
$ws.Range("K42").Value = 22100
$ws.Range("L42").Value = 33900
$ws.Range("M42").Value = 83700
$ws.Range("D43").Value = 1032200
$ws.Range("E43").Value = 1060300
$ws.Range("F43").Value = 985700
$ws.Range("G43").Value = 974200
$ws.Range("H43").Value = 1097100
$ws.Range("I43").Value = 1115800
$ws.Range("J43").Value = 895600
$ws.Range("K43").Value = 808500
$ws.Range("L43").Value = 731200
$ws.Range("M43").Value = 730600
$ws.Range("D44").Value = 29800
$ws.Range("E44").Value = 30000
$ws.Range("F44").Value = 28500
$ws.Range("G44").Value = 26700
$ws.Range("H44").Value = 28300
$ws.Range("I44").Value = 29900
$ws.Range("J44").Value = 31000
$ws.Range("K44").Value = 31500
$ws.Range("L44").Value = 32400
$ws.Range("M44").Value = 34500
$ws.Range("D45").Value = 10500
$ws.Range("E45").Value = 27400
$ws.Range("F45").Value = 19900
$ws.Range("G45").Value = 20100
$ws.Range("H45").Value = 25400
$ws.Range("I45").Value = 25600
$ws.Range("J45").Value = 9500
$ws.Range("K45").Value = 9400
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = 10600
$ws.Range("D46").Value = 1712100
$ws.Range("E46").Value = 1442600
$ws.Range("F46").Value = 1556900
$ws.Range("G46").Value = 1540400
$ws.Range("H46").Value = 1461800
$ws.Range("I46").Value = 1544300
$ws.Range("J46").Value = 1315900
$ws.Range("K46").Value = 1106000
$ws.Range("L46").Value = 1051200
$ws.Range("M46").Value = 1211100
$ws.Range("D47").Value = 1576900
$ws.Range("E47").Value = 2296700
$ws.Range("F47").Value = 2178700
$ws.Range("G47").Value = 2087800
$ws.Range("H47").Value = 2179100
$ws.Range("I47").Value = 2122700
$ws.Range("J47").Value = 2046500
$ws.Range("K47").Value = 1961300
$ws.Range("L47").Value = 1862100
$ws.Range("M47").Value = 1567900
$ws.Range("D48").Value = 2779700
$ws.Range("E48").Value = 2754400
$ws.Range("F48").Value = 2677100
$ws.Range("G48").Value = 2606100
$ws.Range("H48").Value = 2521700
$ws.Range("I48").Value = 2477700
$ws.Range("J48").Value = 2397500
$ws.Range("K48").Value = 2252800
$ws.Range("L48").Value = 2235600
$ws.Range("M48").Value = 2475300
$ws.Range("D49").Value = 1545700
$ws.Range("E49").Value = 1688000
$ws.Range("F49").Value = 1659800
$ws.Range("G49").Value = 1657300
$ws.Range("H49").Value = 1654400
$ws.Range("I49").Value = 1682400
$ws.Range("J49").Value = 1672400
$ws.Range("K49").Value = 1609700
$ws.Range("L49").Value = 1602600
$ws.Range("M49").Value = 1592500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("D52").Value = 1597400
$ws.Range("E52").Value = 779200
$ws.Range("F52").Value = 767400
$ws.Range("G52").Value = 753100
$ws.Range("H52").Value = 752600
$ws.Range("I52").Value = 776100
$ws.Range("J52").Value = 746000
$ws.Range("K52").Value = 744600
$ws.Range("L52").Value = 767400
$ws.Range("M52").Value = 719100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("D54").Value = 9211800
$ws.Range("E54").Value = 8960800
$ws.Range("F54").Value = 8840000
$ws.Range("G54").Value = 8644700
$ws.Range("H54").Value = 8502200
$ws.Range("I54").Value = 8603300
$ws.Range("J54").Value = 8178300
$ws.Range("K54").Value = 7674400
$ws.Range("L54").Value = 7518900
$ws.Range("M54").Value = 7565800
$ws.Range("D57").Value = 363900
$ws.Range("E57").Value = 549500
$ws.Range("F57").Value = 428900
$ws.Range("G57").Value = 356200
$ws.Range("H57").Value = 431600
$ws.Range("I57").Value = 557600
$ws.Range("J57").Value = 391700
$ws.Range("K57").Value = 283700
$ws.Range("L57").Value = 311500
$ws.Range("M57").Value = 299300
$ws.Range("D58").Value = 845500
$ws.Range("E58").Value = 847600
$ws.Range("F58").Value = 1016300
$ws.Range("G58").Value = 659300
$ws.Range("H58").Value = 619600
$ws.Range("I58").Value = 840300
$ws.Range("J58").Value = 755800
$ws.Range("K58").Value = 681600
$ws.Range("L58").Value = 645500
$ws.Range("M58").Value = 490800
$ws.Range("D59").Value = 507100
$ws.Range("E59").Value = 471500
$ws.Range("F59").Value = 535300
$ws.Range("G59").Value = 500100
$ws.Range("H59").Value = 515200
$ws.Range("I59").Value = 567600
$ws.Range("J59").Value = 528000
$ws.Range("K59").Value = 475300
$ws.Range("L59").Value = 446200
$ws.Range("M59").Value = 359100
$ws.Range("D60").Value = 1716500
$ws.Range("E60").Value = 1868600
$ws.Range("F60").Value = 1980400
$ws.Range("G60").Value = 1515500
$ws.Range("H60").Value = 1566500
$ws.Range("I60").Value = 1965500
$ws.Range("J60").Value = 1675500
$ws.Range("K60").Value = 1440600
$ws.Range("L60").Value = 1403100
$ws.Range("M60").Value = 1149100
$ws.Range("D61").Value = 2119600
$ws.Range("E61").Value = 1702900
$ws.Range("F61").Value = 1594900
$ws.Range("G61").Value = 2005500
$ws.Range("H61").Value = 1900700
$ws.Range("I61").Value = 1602200
$ws.Range("J61").Value = 1672600
$ws.Range("K61").Value = 1530600
$ws.Range("L61").Value = 1546800
$ws.Range("M61").Value = 1682700
$ws.Range("D62").Value = 1187400
$ws.Range("E62").Value = 1150600
$ws.Range("F62").Value = 1116300
$ws.Range("G62").Value = 1064200
$ws.Range("H62").Value = 1058400
$ws.Range("I62").Value = 1005900
$ws.Range("J62").Value = 903700
$ws.Range("K62").Value = 885000
$ws.Range("L62").Value = 848300
$ws.Range("M62").Value = 861700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("D66").Value = 5101200
$ws.Range("E66").Value = 4800800
$ws.Range("F66").Value = 4770900
$ws.Range("G66").Value = 4664700
$ws.Range("H66").Value = 4603200
$ws.Range("I66").Value = 4639800
$ws.Range("J66").Value = 4320100
$ws.Range("K66").Value = 3919800
$ws.Range("L66").Value = 3863800
$ws.Range("M66").Value = 3774000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("D72").Value = 1881100
$ws.Range("E72").Value = 1916300
$ws.Range("F72").Value = 1821100
$ws.Range("G72").Value = 1727900
$ws.Range("H72").Value = 1641400
$ws.Range("I72").Value = 1692200
$ws.Range("J72").Value = 1582700
$ws.Range("K72").Value = 1547600
$ws.Range("L72").Value = 1445100
$ws.Range("M72").Value = 1813800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("D76").Value = 4110500
$ws.Range("E76").Value = 4160000
$ws.Range("F76").Value = 4069000
$ws.Range("G76").Value = 3980000
$ws.Range("H76").Value = 3899000
$ws.Range("I76").Value = 3963500
$ws.Range("J76").Value = 3858200
$ws.Range("K76").Value = 3754600
$ws.Range("L76").Value = 3655100
$ws.Range("M76").Value = 3791800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("F80").Value = 43281
$ws.Range("G80").Value = 43190
$ws.Range("H80").Value = 43100
$ws.Range("I80").Value = 43008
$ws.Range("J80").Value = 42916
$ws.Range("K80").Value = 42825
$ws.Range("L80").Value = 42735
$ws.Range("M80").Value = 42643
$ws.Range("D81").Value = 95600
$ws.Range("E81").Value = 91000
$ws.Range("F81").Value = 89000
$ws.Range("G81").Value = 85200
$ws.Range("H81").Value = 19000
$ws.Range("I81").Value = 104900
$ws.Range("J81").Value = 35900
$ws.Range("K81").Value = 101800
$ws.Range("L81").Value = -52400
$ws.Range("M81").Value = -21600
$ws.Range("D83").Value = 49900
$ws.Range("E83").Value = 48000
$ws.Range("F83").Value = 48700
$ws.Range("G83").Value = 45400
$ws.Range("H83").Value = 46700
$ws.Range("I83").Value = 46700
$ws.Range("J83").Value = 47200
$ws.Range("K83").Value = 45400
$ws.Range("L83").Value = 43700
$ws.Range("M83").Value = 44500
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("D89").Value = 43400
$ws.Range("E89").Value = 69600
$ws.Range("F89").Value = 176600
$ws.Range("G89").Value = 164400
$ws.Range("H89").Value = -50000
$ws.Range("I89").Value = 38100
$ws.Range("J89").Value = 156000
$ws.Range("K89").Value = 106000
$ws.Range("L89").Value = 39800
$ws.Range("M89").Value = 57800
$ws.Range("D91").Value = -75000
$ws.Range("E91").Value = -91200
$ws.Range("F91").Value = -86500
$ws.Range("G91").Value = -129000
$ws.Range("H91").Value = -115900
$ws.Range("I91").Value = -48500
$ws.Range("J91").Value = -94000
$ws.Range("K91").Value = -90500
$ws.Range("L91").Value = -695800
$ws.Range("M91").Value = -150400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("D94").Value = -173700
$ws.Range("E94").Value = -137800
$ws.Range("F94").Value = -126700
$ws.Range("G94").Value = -112700
$ws.Range("H94").Value = -52800
$ws.Range("I94").Value = -105600
$ws.Range("J94").Value = -144800
$ws.Range("K94").Value = -98900
$ws.Range("L94").Value = -148600
$ws.Range("M94").Value = -149900
$ws.Range("D96").Value = -4000
$ws.Range("E96").Value = -69100
$ws.Range("F96").Value = -3600
$ws.Range("G96").Value = -500
$ws.Range("H96").Value = -77600
$ws.Range("I96").Value = -1500
$ws.Range("J96").Value = -49900
$ws.Range("K96").Value = -800
$ws.Range("L96").Value = -2400
$ws.Range("M96").Value = -1400
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = 409900
$ws.Range("E100").Value = -148800
$ws.Range("F100").Value = -48100
$ws.Range("G100").Value = 116800
$ws.Range("H100").Value = 34800
$ws.Range("I100").Value = 37100
$ws.Range("J100").Value = 111500
$ws.Range("K100").Value = -16200
$ws.Range("L100").Value = 700
$ws.Range("M100").Value = 238300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = 279600
$ws.Range("E102").Value = -217000
$ws.Range("F102").Value = 1700
$ws.Range("G102").Value = 168500
$ws.Range("H102").Value = -68000
$ws.Range("I102").Value = -30400
$ws.Range("J102").Value = 122600
$ws.Range("K102").Value = -9100
$ws.Range("L102").Value = -108100
$ws.Range("M102").Value = 146200
